$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.181.82'
$ws.Range("E2").Value = '  -6.07%  '
$ws.Range("D3").Value = '2.556.80'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.33%  '
$ws.Range("E6").Value = '  -5.94%  '
$ws.Range("E7").Value = '  -3.43%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.552'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.96%  '
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").Value = '2.947.95'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '2.559.48'
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.874'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.97%  '
$ws.Range("D18").Value = '43.220.47'
$ws.Range("E18").Value = '  -6.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '260.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.75%  '
$ws.Range("E24").Value = '  -4.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '29.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.34%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("E30").Value = '  -3.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.77%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0799'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.15%  '
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.20%  '
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("E43").Value = '  -2.48%  '
$ws.Range("D44").Value = '2.065.96'
$ws.Range("E44").Value = '  -2.57%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '85.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.12%  '
$ws.Range("E47").Value = '  +3.01%  '
$ws.Range("D48").Value = '2.804.21'
$ws.Range("E48").Value = '  -2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '104.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.90%  '
